$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that contain the obsolete Facebook-handle entries.
# Delete from the bottom up so row numbers of the earlier row are not affected.
$ws.Range("A123").EntireRow.Delete()
$ws.Range("A119").EntireRow.Delete()

# Reproduce the author's final selection/scroll position.
$ws.Rows("122:122").Select()
$excel.ActiveWindow.ScrollRow = 253
